# Adds the 2024/12/10 data column (CO) to the "合成確率" sheet, mirroring
# the existing per-day columns: a header date in row 1 and one numeric
# reading per row (2-53), with conditional fill/font matching the value
# (>=140 normal, 125-139.9 light blue, <125 yellow) reused from existing
# cells via copy/paste-special of formats (so styles are byte-identical
# to the ones already used in the sheet, rather than newly synthesized).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One template cell per style already present on the sheet.
$templates = @{ 1 = "A1"; 2 = "BM2"; 3 = "BN2" }

$newColumn = 93   # column CO

# Match column width of the other per-day data columns (width 12).
$ws.Columns($newColumn).ColumnWidth = 11.17

$data = @(
    @{ Row = 1; Value = "2024/12/10"; Style = 1; IsText = $true },
    @{ Row = 2; Value = 136.8; Style = 3; IsText = $false },
    @{ Row = 3; Value = 153.7; Style = 1; IsText = $false },
    @{ Row = 4; Value = 133; Style = 3; IsText = $false },
    @{ Row = 5; Value = 149.8; Style = 1; IsText = $false },
    @{ Row = 6; Value = 124.1; Style = 2; IsText = $false },
    @{ Row = 7; Value = 121.5; Style = 2; IsText = $false },
    @{ Row = 8; Value = 197.1; Style = 1; IsText = $false },
    @{ Row = 9; Value = 151; Style = 1; IsText = $false },
    @{ Row = 10; Value = 224.8; Style = 1; IsText = $false },
    @{ Row = 11; Value = 178.6; Style = 1; IsText = $false },
    @{ Row = 12; Value = 182.3; Style = 1; IsText = $false },
    @{ Row = 13; Value = 149.7; Style = 1; IsText = $false },
    @{ Row = 14; Value = 184.4; Style = 1; IsText = $false },
    @{ Row = 15; Value = 136.8; Style = 3; IsText = $false },
    @{ Row = 16; Value = 160.1; Style = 1; IsText = $false },
    @{ Row = 17; Value = 179.8; Style = 1; IsText = $false },
    @{ Row = 18; Value = 220.8; Style = 1; IsText = $false },
    @{ Row = 19; Value = 153.8; Style = 1; IsText = $false },
    @{ Row = 20; Value = 160.7; Style = 1; IsText = $false },
    @{ Row = 21; Value = 136.7; Style = 3; IsText = $false },
    @{ Row = 22; Value = 149; Style = 1; IsText = $false },
    @{ Row = 23; Value = 175.6; Style = 1; IsText = $false },
    @{ Row = 24; Value = 174.4; Style = 1; IsText = $false },
    @{ Row = 25; Value = 178.5; Style = 1; IsText = $false },
    @{ Row = 26; Value = 204.2; Style = 1; IsText = $false },
    @{ Row = 27; Value = 136.9; Style = 3; IsText = $false },
    @{ Row = 28; Value = 154.3; Style = 1; IsText = $false },
    @{ Row = 29; Value = 168; Style = 1; IsText = $false },
    @{ Row = 30; Value = 132.1; Style = 3; IsText = $false },
    @{ Row = 31; Value = 155.6; Style = 1; IsText = $false },
    @{ Row = 32; Value = 170.5; Style = 1; IsText = $false },
    @{ Row = 33; Value = 143.8; Style = 1; IsText = $false },
    @{ Row = 34; Value = 142.6; Style = 1; IsText = $false },
    @{ Row = 35; Value = 119.4; Style = 2; IsText = $false },
    @{ Row = 36; Value = 171.5; Style = 1; IsText = $false },
    @{ Row = 37; Value = 192.7; Style = 1; IsText = $false },
    @{ Row = 38; Value = 193; Style = 1; IsText = $false },
    @{ Row = 39; Value = 340.6; Style = 1; IsText = $false },
    @{ Row = 40; Value = 215.4; Style = 1; IsText = $false },
    @{ Row = 41; Value = 130; Style = 3; IsText = $false },
    @{ Row = 42; Value = 137.7; Style = 3; IsText = $false },
    @{ Row = 43; Value = 153.1; Style = 1; IsText = $false },
    @{ Row = 44; Value = 140.6; Style = 1; IsText = $false },
    @{ Row = 45; Value = 178.3; Style = 1; IsText = $false },
    @{ Row = 46; Value = 131.3; Style = 3; IsText = $false },
    @{ Row = 47; Value = 110.4; Style = 2; IsText = $false },
    @{ Row = 48; Value = 114.4; Style = 2; IsText = $false },
    @{ Row = 49; Value = 181.8; Style = 1; IsText = $false },
    @{ Row = 50; Value = 163.1; Style = 1; IsText = $false },
    @{ Row = 51; Value = 154.4; Style = 1; IsText = $false },
    @{ Row = 52; Value = 222.2; Style = 1; IsText = $false },
    @{ Row = 53; Value = 136.7; Style = 3; IsText = $false }
)

foreach ($item in $data) {
    $dst = $ws.Range("CO" + $item.Row)
    $src = $ws.Range($templates[$item.Style])

    if ($item.IsText) {
        # Leading apostrophe forces text entry so the date-like string
        # isn't auto-converted into a date serial number.
        $dst.Value = "'" + $item.Value
        $src.Copy()
        $dst.PasteSpecial(-4122)   # xlPasteFormats
    } else {
        $src.Copy()
        $dst.PasteSpecial(-4122)   # xlPasteFormats
        $dst.Value = $item.Value
    }
}

$excel.CutCopyMode = $false
